$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function InsertRowWithFormat($targetRow, $templateRow) {
    $ws.Rows("$targetRow`:$targetRow").Insert()
    $ws.Range("A$templateRow`:H$templateRow").Copy()
    $ws.Range("A$targetRow`:H$targetRow").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# --- Header / summary updates ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 22

# --- Zero-out existing Wednesday table pricing values ---
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0

# --- Insert 5 new rows into the Wednesday table, just above its TOTAL row (21) ---
InsertRowWithFormat 21 19
InsertRowWithFormat 22 20
InsertRowWithFormat 23 21
InsertRowWithFormat 24 22
InsertRowWithFormat 25 23

# --- Insert 4 new rows into the Friday table, just above its TOTAL row (39 after the shift above) ---
InsertRowWithFormat 39 37
InsertRowWithFormat 40 38
InsertRowWithFormat 41 39
InsertRowWithFormat 42 40

# --- Fill in the new Wednesday rows ---
$ws.Range("A21").Value = "Point 15"
$ws.Range("B21").Value = "ANC-DHM-8-84-T1-C"
$ws.Range("C21").Value = "Inst"
$ws.Range("D21").Value = "ANC,Dbl Hlx Mach,8in,84in,TpEye 1in,Cor"
$ws.Range("E21").Value = "EA"
$ws.Range("F21").Value = 0
$ws.Range("H21").Value = 0

$ws.Range("A22").Value = "Point 31"
$ws.Range("B22").Value = "PIN-15-PTP-C"
$ws.Range("C22").Value = "Inst"
$ws.Range("D22").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("E22").Value = "EA"
$ws.Range("F22").Value = 0
$ws.Range("H22").Value = 0

$ws.Range("A23").Value = "Point 33"
$ws.Range("B23").Value = "PIN-15-PTP-C"
$ws.Range("C23").Value = "Inst"
$ws.Range("D23").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("E23").Value = "EA"
$ws.Range("F23").Value = 0
$ws.Range("H23").Value = 0

$ws.Range("A24").Value = "Point 01"
$ws.Range("B24").Value = "GYA-HDIG"
$ws.Range("C24").Value = "Inst"
$ws.Range("D24").Value = "GYA, Hand Dig or Additional Excavation"
$ws.Range("E24").Value = "EA"
$ws.Range("F24").Value = 1
$ws.Range("H24").Value = 0

$ws.Range("A25").Value = "Point 03"
$ws.Range("B25").Value = "PIN-15-PTP-C"
$ws.Range("C25").Value = "Inst"
$ws.Range("D25").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("E25").Value = "EA"
$ws.Range("F25").Value = 0
$ws.Range("H25").Value = 0

# --- Wednesday TOTAL row (now at 26) ---
$ws.Range("H26").Value = 0

# --- Zero-out existing Friday table pricing values (rows 31-38 after shift) ---
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("H38").Value = 0

# --- Fill in the new Friday rows ---
$ws.Range("A39").Value = "Point 35"
$ws.Range("B39").Value = "INS-15-D-S-C"
$ws.Range("C39").Value = "Inst"
$ws.Range("D39").Value = "INS,15kV,Deadend,Polymer,Corr"
$ws.Range("E39").Value = "EA"
$ws.Range("F39").Value = 1
$ws.Range("H39").Value = 0

$ws.Range("A40").Value = "Point 35"
$ws.Range("B40").Value = "POL-40-2"
$ws.Range("C40").Value = "Inst"
$ws.Range("D40").Value = "Pole,40ft,Class 2"
$ws.Range("E40").Value = "EA"
$ws.Range("F40").Value = 0
$ws.Range("H40").Value = 0

$ws.Range("A41").Value = "Point 17"
$ws.Range("B41").Value = "GYA-HDIG"
$ws.Range("C41").Value = "Inst"
$ws.Range("D41").Value = "GYA, Hand Dig or Additional Excavation"
$ws.Range("E41").Value = "EA"
$ws.Range("F41").Value = 1
$ws.Range("H41").Value = 0

$ws.Range("A42").Value = "Point 03"
$ws.Range("B42").Value = "GYA-HDIG"
$ws.Range("C42").Value = "Inst"
$ws.Range("D42").Value = "GYA, Hand Dig or Additional Excavation"
$ws.Range("E42").Value = "EA"
$ws.Range("F42").Value = 1
$ws.Range("H42").Value = 0

# --- Friday TOTAL row (now at 43) ---
$ws.Range("H43").Value = 0
